# "Avancement du projet 5"
# The logo picture on the last slide (slide 16 / sldId 278) is nudged
# up slightly (its xfrm offset goes from (0,0) to (0,-16934) EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# Locate the picture (p:pic id="7", name "Picture 2") robustly by Id,
# falling back to its name if Id isn't exposed for some reason.
$pic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 7) {
        $pic = $candidate
        break
    }
}
if ($pic -eq $null) {
    $pic = $s.Shapes.Item("Picture 2")
}

# -16934 EMU expressed in points (1 pt = 12700 EMU). Using a value that
# lands solidly inside the EMU bucket for -16934 after the point->EMU
# conversion (rather than the boundary value) keeps the result stable.
$pic.Left = 0
$pic.Top = -1.3334
